# riska.xlsx update: refresh the repayment collector stats table and bump
# the workbook's "(2)" re-export counter to "(4)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric cell updates -------------------------------------------------
# Repayment_collections (D), Talk_time (H) and Repayment_new_collections (J)
# are plain numbers in the sheet, so a normal .Value assignment is correct.
$numericUpdates = @{
    "H2"  = 660
    "D3"  = 5
    "H3"  = 488
    "D4"  = 9
    "H4"  = 741
    "D5"  = 6
    "H5"  = 381
    "D6"  = 7
    "H6"  = 511
    "D7"  = 8
    "H7"  = 725
    "J7"  = 1
    "D8"  = 10
    "H8"  = 1.135
    "H9"  = 1.107
    "D10" = 8
    "H10" = 320
    "H11" = 251
    "H12" = 357
    "J12" = 1
    "H13" = 550
    "H14" = 767
    "H15" = 2.2320000000000002
    "H16" = 973
    "D17" = 7
    "H17" = 437
    "D18" = 2
    "H18" = 870
    "J18" = 1
}

foreach ($addr in $numericUpdates.Keys) {
    $ws.Range($addr).Value = $numericUpdates[$addr]
}

# --- Text cell updates -----------------------------------------------------
# Repayment_amount (E), Pending Amount Recovery (G), New_collection_amount_rate
# (K) and New_collection_count_rate (L) are stored as literal text (e.g.
# "3,309,416.00") in this workbook rather than real numbers, so force the
# Text number format before assigning the value - otherwise Excel helpfully
# (and, here, unhelpfully) reinterprets the comma/decimal string as a number.
# The style is reset back to Normal afterwards so no stray formatting sticks
# to the cell.
$textUpdates = @{
    "E3"  = "3,309,416.00"
    "G3"  = "1.84"
    "E4"  = "1,113,098.00"
    "G4"  = "0.81"
    "E5"  = "1,696,943.00"
    "G5"  = "1.32"
    "E6"  = "1,537,622.00"
    "G6"  = "1.18"
    "E7"  = "1,461,376.00"
    "G7"  = "1.23"
    "K7"  = "3.35"
    "L7"  = "7.14"
    "E8"  = "4,940,492.00"
    "G8"  = "2.80"
    "K9"  = "3.02"
    "E10" = "1,440,650.00"
    "G10" = "0.99"
    "K12" = "2.17"
    "L12" = "6.67"
    "E17" = "2,061,127.00"
    "G17" = "1.34"
    "E18" = "653,237.00"
    "G18" = "0.49"
    "L18" = "7.14"
}

foreach ($addr in $textUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$addr]
    $cell.Style = "Normal"
}

# --- Workbook/sheet rename --------------------------------------------------
# The re-uploaded file bumped the "(2)" suffix in both the sheet tab name and
# (implicitly, via Excel) the workbook's title-of-parts metadata.
$ws.Name = "repayment_20250926_20250926 (4)"
